# Generate Report for Archive
# The localization status for the file "7f7468c8-c24a-4952-b815-cfa0687a9e07.md"
# (row 5 in every sheet) moved from "Ready for handoff" to "In Translation".
# Update the Status column on the per-locale sheets (zh-cn, de-de) and the
# corresponding roll-up columns on the Overview sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E5").Value = "In Translation"
$overview.Range("F5").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C5").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C5").Value = "In Translation"
